$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every Price (D) and Volume(1h) (E) figure as
# literal text (coinranking.com scrape formatting, e.g. thousands-dot
# separators and padded "  +x.xx%  " strings) even when a cell happens to
# look like a plain decimal number. Force both columns to Text format up
# front so Excel does not silently re-interpret values such as "247.88"
# or "0.100" as numbers (which would also eat significant trailing zeros).
$ws.Range('D2:D51').NumberFormat = '@'
$ws.Range('E2:E51').NumberFormat = '@'

# --- Update Price (D) and Volume(1h) (E) columns for most rows ---
$ws.Range('D2').Value = '36.893.84'
$ws.Range('E2').Value = '  +4.32%  '
$ws.Range('D3').Value = '1.913.03'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '247.88'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').Value = '0.683'
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '47.36'
$ws.Range('E8').Value = '  +9.65%  '
$ws.Range('D9').Value = '0.373'
$ws.Range('E9').Value = '  +5.17%  '
$ws.Range('D10').Value = '58.13'
$ws.Range('E10').Value = '  +6.07%  '
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('D12').Value = '0.100'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').Value = '15.53'
$ws.Range('E13').Value = '  +12.88%  '
$ws.Range('D14').Value = '0.817'
$ws.Range('E14').Value = '  +6.24%  '
$ws.Range('D15').Value = '2.190.71'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '5.09'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '1.920.27'
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('D18').Value = '36.821.35'
$ws.Range('E18').Value = '  +4.18%  '
$ws.Range('D19').Value = '74.33'
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('E21').Value = '  +5.78%  '
$ws.Range('D22').Value = '249.45'
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '2.47'
$ws.Range('E25').Value = '  -6.43%  '
$ws.Range('D26').Value = '167.02'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('E27').Value = '  -3.15%  '
$ws.Range('D28').Value = '8.75'
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '18.59'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  +5.77%  '
$ws.Range('D32').Value = '0.0606'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').Value = '0.0909'
$ws.Range('E33').Value = '  +26.78%  '
$ws.Range('D34').Value = '4.25'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '18.75'
$ws.Range('E37').Value = '  +36.52%  '
$ws.Range('D38').Value = '0.876'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('D43').Value = '17.43'
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('D44').Value = '2.88'
$ws.Range('E44').Value = '  +19.93%  '
$ws.Range('D45').Value = '1.08'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').Value = '1.345.73'
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').Value = '2.40'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('D48').Value = '0.0832'
$ws.Range('E48').Value = '  +2.83%  '
$ws.Range('D49').Value = '2.79'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('D50').Value = '6.35'
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('D51').Value = '2.102.74'
$ws.Range('E51').Value = '  +2.21%  '

# --- Rows 40/41: Aave and LidoDAOToken swap ranking positions ---
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '106.09'
$ws.Range('E40').Value = '  +8.68%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '1.94'
$ws.Range('E41').Value = '  +0.47%  '
